$d = $word.ActiveDocument

$replacements = @(
    @{old="22÷8=2, 6"; new="30÷4=7, 2"},
    @{old="83÷2=41, 1"; new="86÷5=17, 1"},
    @{old="75÷6=12, 3"; new="12÷6=2, 0"},
    @{old="56÷9=6, 2"; new="73÷9=8, 1"},
    @{old="10÷3=3, 1"; new="75÷4=18, 3"},
    @{old="54÷3=18, 0"; new="83÷7=11, 6"},
    @{old="19÷5=3, 4"; new="13÷4=3, 1"},
    @{old="59÷6=9, 5"; new="62÷4=15, 2"},
    @{old="55÷4=13, 3"; new="42÷4=10, 2"},
    @{old="16÷2=8, 0"; new="78÷6=13, 0"},
    @{old="69÷3=23, 0"; new="43÷7=6, 1"},
    @{old="63÷2=31, 1"; new="13÷8=1, 5"},
    @{old="30÷3=10, 0"; new="21÷8=2, 5"},
    @{old="25÷5=5, 0"; new="16÷4=4, 0"},
    @{old="86÷7=12, 2"; new="95÷3=31, 2"},
    @{old="96÷5=19, 1"; new="80÷9=8, 8"},
    @{old="97÷5=19, 2"; new="65÷8=8, 1"},
    @{old="64÷6=10, 4"; new="58÷7=8, 2"},
    @{old="91÷8=11, 3"; new="43÷4=10, 3"},
    @{old="46÷3=15, 1"; new="41÷6=6, 5"},
    @{old="55÷9=6, 1"; new="45÷5=9, 0"},
    @{old="89÷5=17, 4"; new="26÷9=2, 8"},
    @{old="50÷8=6, 2"; new="34÷4=8, 2"},
    @{old="73÷5=14, 3"; new="66÷9=7, 3"},
    @{old="37÷4=9, 1"; new="61÷7=8, 5"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
